# updated the dates and links in the presentations
$p = $ppt.ActivePresentation

# --- Slide 1 ("Intro to R for Biologists / Session 2") ---
# Shape 2 ("CustomShape 2") holds the presenter names / term-year line.
# Update "Hilary 2022" -> "Trinity 2022" while preserving the run's
# existing character formatting (color/typeface/size).
$s1 = $p.Slides.Item(1)
$shp1b = $s1.Shapes.Item(2)
$tr1b = $shp1b.TextFrame.TextRange
$full1b = $tr1b.Text
$oldTerm = "Hilary 2022"
$newTerm = "Trinity 2022"
$idx = $full1b.IndexOf($oldTerm)
if ($idx -ge 0) {
    $termRange = $tr1b.Characters($idx + 1, $oldTerm.Length)
    $termRange.Text = $newTerm
}

# --- Slide 20 ("Resources") ---
# Shape 2 ("CustomShape 2") lists resource links. The trailing
# "(registration open 2nd Feb)" annotation (including the leading
# space that separated it from the course title) after "Advanced R
# course for data analysis and visualisation" has been removed.
$s20 = $p.Slides.Item(20)
$shp20b = $s20.Shapes.Item(2)
$tr20b = $shp20b.TextFrame.TextRange
$fullText20 = $tr20b.Text
$anchor = "Advanced R course for data analysis and visualisation"
$trailing = " (registration open 2nd Feb)"
$anchorIdx = $fullText20.IndexOf($anchor)
if ($anchorIdx -ge 0) {
    $afterIdx = $anchorIdx + $anchor.Length
    $remainder = $fullText20.Substring($afterIdx)
    if ($remainder.StartsWith($trailing)) {
        $removeRange = $tr20b.Characters($afterIdx + 1, $trailing.Length)
        $removeRange.Text = ""
    }
}
